# Append: 2025-12-09 12:52 JST
# Update the "取得日時" (acquired timestamp) column on the ランサーズ sheet
# from 2025-12-09 12:39:44 to 2025-12-09 12:52:44 for all existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-09 12:39:44"
$newTimestamp = "2025-12-09 12:52:44"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
